# Apply the update described in the commit: insert a new "LB" column
# after "safe_time" (before "obj"), shifting obj/CI/train_time/test_time
# one column to the right, and refresh all the numeric results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing obj/CI/train_time/test_time headers one column to
# the right (D->E, E->F, F->G, G->H) before writing the new "LB" header
# into column D. Work from the rightmost column backwards so values are
# not clobbered before they are copied. The numeric data underneath the
# headers is rewritten explicitly afterwards with refreshed results.
$ws.Range("H1").Value = $ws.Range("G1").Value2
$ws.Range("G1").Value = $ws.Range("F1").Value2
$ws.Range("F1").Value = $ws.Range("E1").Value2
$ws.Range("E1").Value = $ws.Range("D1").Value2
$ws.Range("D1").Value = "LB"

# New LB column values (rows 2-9)
$lb = @(
    2886606.1687525799,
    3279854.349891,
    3600066.4793544598,
    3858265.7698462298,
    2663182.0271669799,
    3124531.7976661902,
    3439902.7263436099,
    3690635.1377503299
)

for ($i = 0; $i -lt $lb.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $lb[$i]
}

# Updated obj / CI / train_time / test_time values (now columns E, F, G, H)
$data = @(
    @(2696304.9390490302, 161299.519050559,  139.355339050292,   6.7359700202941797),
    @(3100418.8371999399, 183713.694034418,  99.673047065734806, 6.2072639465331996),
    @(3466913.4481158601, 206524.58693778201,92.225600957870398, 6.3155429363250697),
    @(3667496.0065426598, 202221.74610881301,100.651798963546,   6.2693567276000897),
    @(2568513.1194338398, 223768.42324649799,107.20077967643699, 8.1055607795715297),
    @(3078117.6668686401, 263259.74107469001,61.789297342300401, 6.77239966392517),
    @(3293880.6603405899, 216223.50241273499,26.608678340911801, 5.5276601314544598),
    @(3494904.0254636202, 206203.06557664901,27.719955921173,    5.5746469497680602)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $data[$i][0]
    $ws.Cells.Item($row, 6).Value = $data[$i][1]
    $ws.Cells.Item($row, 7).Value = $data[$i][2]
    $ws.Cells.Item($row, 8).Value = $data[$i][3]
}

# Update the selection to D2 to match the saved workbook view
$ws.Range("D2").Select()

$wb.Save()
